# according to board: pho-klite, change the pinout
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: PB1 -> PA4, the "一体板" light-sensor description changes to the "外部" one,
# and column E (丝印/silkscreen) now gets an "ADC2" label.
$ws.Range("B4").Value = "PA4"
$ws.Range("D4").Value = "外部光敏电阻输入,ADC12_IN4"
$ws.Range("E4").Value = "ADC2"

# New row 7: beeper output on PA0
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "PA0"
$ws.Range("C7").Value = "输出"
$ws.Range("D7").Value = "beeper"
$ws.Range("E7").Value = "ADC1"

# Copy the bordered style of the existing rows down onto the new row
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Column E width adjustment (became custom, wider, to fit "ADC2"/"ADC1" labels)
$ws.Range("E1").EntireColumn.ColumnWidth = 9.86

# Update the selected cell as recorded in the workbook view
$ws.Range("D13").Select()
